# Updated Master database. Adjusted update-database script to position
# sample names in correct sample description columns.

$wb = $excel.ActiveWorkbook

# --- client_list / site_codes: add new clients + their site codes ---------
# New string values are interleaved in the same order the original author
# entered them (new clients, then their codes) so that the shared-string
# table is built up in the same sequence as the source workbook.
$clientList = $wb.Worksheets.Item("client_list")
[void]$clientList.Activate()

$siteCodes = $wb.Worksheets.Item("site_codes")

# Fix "Lanscape" -> "Landscape" typo on the existing Town of Basalt row.
$siteCodes.Cells.Item(47, 3).Value = "Landscape"

$clientList.Cells.Item(42, 1).Value = 85137
$clientList.Cells.Item(42, 2).Value = "Evergro"

$clientList.Cells.Item(43, 1).Value = 93470
$clientList.Cells.Item(43, 2).Value = "Phillips Lawn Sprinkler Co"

$siteCodes.Cells.Item(48, 1).Value = "Evergro"
$siteCodes.Cells.Item(48, 2).Value = "EVERGRO"
$siteCodes.Cells.Item(48, 3).Value = "Landscape"

$siteCodes.Cells.Item(49, 1).Value = "Phillips Lawn Sprinkler Co"
$siteCodes.Cells.Item(49, 2).Value = "PHILLIPS"
$siteCodes.Cells.Item(49, 3).Value = "Landscape"

$clientList.Cells.Item(44, 1).Value = 74400
$clientList.Cells.Item(44, 2).Value = "Living Earth Dallas"

$siteCodes.Cells.Item(50, 1).Value = "Living Earth Dallas"
$siteCodes.Cells.Item(50, 2).Value = "LEDALLAS"
$siteCodes.Cells.Item(50, 3).Value = "Landscape"

[void]$clientList.Range("B40").Select()

[void]$siteCodes.Activate()
[void]$siteCodes.Range("A57").Select()

# --- catalog: becomes the active/selected sheet in the saved file ---------
$catalog = $wb.Worksheets.Item("catalog")
[void]$catalog.Activate()
[void]$catalog.Range("C621").Select()
